$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# New weekly WRESBAL observations to append to the "Data" sheet.
$newRows = @(
    @{ Row = 104; Date = 45189; Value = 3231.649 },
    @{ Row = 105; Date = 45196; Value = 3170.324 },
    @{ Row = 106; Date = 45203; Value = 3145.72 },
    @{ Row = 107; Date = 45210; Value = 3288.945 },
    @{ Row = 108; Date = 45217; Value = 3353.881 },
    @{ Row = 109; Date = 45224; Value = 3261.886 }
)

# Use the formatting already applied to the last existing data row (A103) as
# the template for the new date cells so the new rows look identical in style.
$dateStyleSource = $ws.Range("A103")

foreach ($item in $newRows) {
    $r = $item.Row
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    # Assign Value2 (not Value) so Excel does not try to auto-detect/auto
    # format the plain date-serial number while we are still applying the
    # explicit format below.
    $aCell.Value2 = $item.Date

    $dateStyleSource.Copy()
    $aCell.PasteSpecial(-4122)

    $bCell.Value2 = $item.Value
}

$excel.CutCopyMode = 0

# Update SeriesInfo sheet metadata to reflect the refreshed pull. The
# target values look like dates/timestamps, and assigning them directly
# to .Value lets Excel "helpfully" reinterpret them as date serials with
# date formatting. To keep them as plain text (matching the original
# cells' storage), build each value through a text formula and then
# flatten the formula down to its static value via copy/paste-values.
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

Set-TextValue $wsInfo.Range("B3") "2023-10-27"
Set-TextValue $wsInfo.Range("B4") "2023-10-27"
Set-TextValue $wsInfo.Range("B7") "2023-10-25"
Set-TextValue $wsInfo.Range("B14") "2023-10-26 15:35:02-05"

$wsInfo.Range("B15").Value2 = 73
